$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.032.48"
$ws.Range("E2").Value = "  -2.47%  "

$ws.Range("D3").Value = "2.504.52"
$ws.Range("E3").Value = "  -3.87%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'552.50"
$ws.Range("E5").Value = "  -3.49%  "

$ws.Range("D6").Value = "'147.92"
$ws.Range("E6").Value = "  -4.77%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -3.74%  "

$ws.Range("D9").Value = "2.502.87"
$ws.Range("E9").Value = "  -3.86%  "

$ws.Range("D10").Value = "'0.108"
$ws.Range("E10").Value = "  -8.63%  "

$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").Value = "'5.41"
$ws.Range("E12").Value = "  -7.27%  "

$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  -6.62%  "

$ws.Range("D14").Value = "'26.28"
$ws.Range("E14").Value = "  -6.65%  "

$ws.Range("D15").Value = "2.951.34"
$ws.Range("E15").Value = "  -4.22%  "

$ws.Range("D16").Value = "61.878.38"
$ws.Range("E16").Value = "  -2.50%  "

$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -7.67%  "

$ws.Range("D18").Value = "2.484.61"
$ws.Range("E18").Value = "  -4.98%  "

$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  -6.34%  "

$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "  -6.48%  "

$ws.Range("D21").Value = "'4.21"
$ws.Range("E21").Value = "  -7.54%  "

$ws.Range("D22").Value = "'322.71"
$ws.Range("E22").Value = "  -5.82%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'63.82"
$ws.Range("E24").Value = "  -5.06%  "

$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  -1.21%  "

$ws.Range("D26").Value = "'0.0000104"
$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("D27").Value = "2.605.43"
$ws.Range("E27").Value = "  -5.69%  "

$ws.Range("D28").Value = "'1.51"
$ws.Range("E28").Value = "  -3.41%  "

$ws.Range("D29").Value = "'543.24"
$ws.Range("E29").Value = "  -7.53%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'8.45"
$ws.Range("E30").Value = "  -7.61%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").Value = "'7.69"
$ws.Range("E32").Value = "  -2.55%  "

$ws.Range("D33").Value = "'0.148"
$ws.Range("E33").Value = "  -7.96%  "

$ws.Range("E34").Value = "  -7.88%  "

$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  -7.69%  "

$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  -8.87%  "

$ws.Range("D37").Value = "'4.92"
$ws.Range("E37").Value = "  -9.28%  "

$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Value = "'0.379"
$ws.Range("E39").Value = "  -6.31%  "

$ws.Range("D40").Value = "'18.57"
$ws.Range("E40").Value = "  -5.66%  "

$ws.Range("D41").Value = "'144.58"
$ws.Range("E41").Value = "  -6.94%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.72"
$ws.Range("E42").Value = "  -7.58%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'40.59"
$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").Value = "'2.36"
$ws.Range("E45").Value = "  -5.39%  "

$ws.Range("D46").Value = "'149.52"
$ws.Range("E46").Value = "  -4.22%  "

$ws.Range("D47").Value = "'3.59"
$ws.Range("E47").Value = "  -7.99%  "

$ws.Range("D48").Value = "'21.26"
$ws.Range("E48").Value = "  -7.90%  "

$ws.Range("D49").Value = "'0.0539"
$ws.Range("E49").Value = "  -8.29%  "

$ws.Range("D50").Value = "'0.592"
$ws.Range("E50").Value = "  -5.73%  "

$ws.Range("D51").Value = "'0.0947"
$ws.Range("E51").Value = "  -5.53%  "
